$d = $word.ActiveDocument

# --- Edit 1: "Url:" paragraph - remove pl-PL lang formatting ---
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:pPr><w:spacing w:before="200"/></w:pPr><w:r><w:t xml:space="preserve">Url: </w:t></w:r><w:hyperlink r:id="rId5" w:history="1"><w:r><w:t>http://demo1a3.project.experienceit.pl/</w:t></w:r></w:hyperlink></w:p>'
$p1 = $d.Paragraphs(5)
$p1.Range.InsertXML($xml1)

# --- Edit 2: "I've got this token" paragraph - merge runs, remove pl-PL lang/proofErr ---
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="200"/></w:pPr><w:r><w:t>I’ve got this token:</w:t></w:r><w:r><w:t xml:space="preserve"> JKN0jaysntR3ln_Wwxi7e_jA_9D9YtaY6pLH0TW90Lo</w:t></w:r></w:p>'
$p2 = $d.Paragraphs(10)
$p2.Range.InsertXML($xml2)

# --- Edit 3: replace "[TODO]" paragraph with Pathauto instructions (14 paragraphs) ---
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="200"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Drupal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 8 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RESTful</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Web Services</w:t></w:r><w:r><w:t xml:space="preserve"> doesn’t allow to change path (permalink) while creating new content. So we need to do that after content is created. We will use </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pathauto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> module to do that. </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="200"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Inside this project there is “module” folder that contains </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pathauto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> module along with two modules required by </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pathauto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ctools</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and token. These three modules are in my every </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Drupal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> installation, so you would probably install them anyway. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Nagwek2"/><w:spacing w:before="360" w:after="280"/></w:pPr><w:r><w:t xml:space="preserve">6.1 Enable </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pathauto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="200"/></w:pPr><w:r><w:t xml:space="preserve">Copy whole “module” folder to your </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>drupal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> installation.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="200"/></w:pPr><w:r><w:t xml:space="preserve">Go to: </w:t></w:r><w:r><w:t>/admin/modules</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="200"/></w:pPr><w:r><w:t xml:space="preserve">Enable </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Pathauto</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="200"/></w:pPr><w:r><w:t>Click “Install”</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Nagwek2"/><w:spacing w:before="360" w:after="280"/></w:pPr><w:r><w:t xml:space="preserve">6.2 Configure </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pathauto</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="200"/></w:pPr><w:r><w:t xml:space="preserve">Go to: </w:t></w:r><w:r><w:t>/admin/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>config</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/search/path/patterns/add</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="200"/></w:pPr><w:r><w:t>Pattern type: Content</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="200"/></w:pPr><w:r><w:t xml:space="preserve">Path pattern: </w:t></w:r><w:r><w:t>invoices/[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>node:field_no</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>]</w:t></w:r><w:r><w:t xml:space="preserve"> – </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>field_no</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is invoice number. So if your invoice number is for example “inv 005” - </w:t></w:r><w:r><w:t xml:space="preserve">for every new invoice this will create permalink </w:t></w:r><w:r><w:t>like invoices/inv-005. Feel free to use other fields from this node.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="200"/></w:pPr><w:r><w:t>Content type – select “Invoice”</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="200"/></w:pPr><w:r><w:t>Label: “</w:t></w:r><w:r><w:t>Invoices</w:t></w:r><w:r><w:t>”</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="200"/></w:pPr><w:r><w:t>Click “Save”</w:t></w:r></w:p>'
$p3 = $d.Paragraphs(141)
$p3.Range.InsertXML($xml3)

# --- Edit 4: "[INCOMPLETE]" paragraph - drop lastRenderedPageBreak, add spacing ---
$xml4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="200"/></w:pPr><w:r><w:t>[INCOMPLETE] Inside this project there is “template” folder - with invoice template inside. Copy whole “template” folder into front-end theme. I’m using “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Bartik</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” theme so I would copy “template” folder to /core/themes/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bartik</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. There is already “template” folder, so invoice template will be added to other </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bartik</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> templates.</w:t></w:r></w:p>'
$p4 = $d.Paragraphs(143 + 13)
$p4.Range.InsertXML($xml4)

Write-Host "Edits applied"
